# Update rows 2-7 with new values and add new rows 8-13,
# reflecting the expanded Gnai2-Cnr1 ligand-receptor pair analysis
# (Natmi re-run following Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "Cnr1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 106.8785896666667
$ws.Cells.Item(2, 8).Value = 320.635769
$ws.Cells.Item(2, 9).Value = 0.1508748302900445
$ws.Cells.Item(2, 10).Value = 0.1508748302900445
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.03992133333333333
$ws.Cells.Item(2, 14).Value = 0.119764
$ws.Cells.Item(2, 15).Value = 0.01377024262113379
$ws.Cells.Item(2, 16).Value = 0.01377024262113379
$ws.Cells.Item(2, 17).Value = 4.266735804279555
$ws.Cells.Item(2, 18).Value = 38.40062223851599
$ws.Cells.Item(2, 19).Value = 0.002077583018516298
$ws.Cells.Item(2, 20).Value = 0.002077583018516298

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "Cnr1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 106.8785896666667
$ws.Cells.Item(3, 8).Value = 320.635769
$ws.Cells.Item(3, 9).Value = 0.1508748302900445
$ws.Cells.Item(3, 10).Value = 0.1508748302900445
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.859180333333333
$ws.Cells.Item(3, 14).Value = 8.577541
$ws.Cells.Item(3, 15).Value = 0.9862297573788662
$ws.Cells.Item(3, 16).Value = 0.9862297573788662
$ws.Cells.Item(3, 17).Value = 305.5851616293365
$ws.Cells.Item(3, 18).Value = 2750.266454664029
$ws.Cells.Item(3, 19).Value = 0.1487972472715282
$ws.Cells.Item(3, 20).Value = 0.1487972472715282

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "Cnr1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 77.232923
$ws.Cells.Item(4, 8).Value = 231.698769
$ws.Cells.Item(4, 9).Value = 0.1090256166999485
$ws.Cells.Item(4, 10).Value = 0.1090256166999485
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.03992133333333333
$ws.Cells.Item(4, 14).Value = 0.119764
$ws.Cells.Item(4, 15).Value = 0.01377024262113379
$ws.Cells.Item(4, 16).Value = 0.01377024262113379
$ws.Cells.Item(4, 17).Value = 3.083241263390666
$ws.Cells.Item(4, 18).Value = 27.749171370516
$ws.Cells.Item(4, 19).Value = 0.001501309193877027
$ws.Cells.Item(4, 20).Value = 0.001501309193877027

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "Cnr1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 77.232923
$ws.Cells.Item(5, 8).Value = 231.698769
$ws.Cells.Item(5, 9).Value = 0.1090256166999485
$ws.Cells.Item(5, 10).Value = 0.1090256166999485
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.859180333333333
$ws.Cells.Item(5, 14).Value = 8.577541
$ws.Cells.Item(5, 15).Value = 0.9862297573788662
$ws.Cells.Item(5, 16).Value = 0.9862297573788662
$ws.Cells.Item(5, 17).Value = 220.8228545274476
$ws.Cells.Item(5, 18).Value = 1987.405690747029
$ws.Cells.Item(5, 19).Value = 0.1075243075060715
$ws.Cells.Item(5, 20).Value = 0.1075243075060715

# Row 6
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "Cnr1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 136.676337
$ws.Cells.Item(6, 8).Value = 410.029011
$ws.Cells.Item(6, 9).Value = 0.1929387280825172
$ws.Cells.Item(6, 10).Value = 0.1929387280825172
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.03992133333333333
$ws.Cells.Item(6, 14).Value = 0.119764
$ws.Cells.Item(6, 15).Value = 0.01377024262113379
$ws.Cells.Item(6, 16).Value = 0.01377024262113379
$ws.Cells.Item(6, 17).Value = 5.456301608156
$ws.Cells.Item(6, 18).Value = 49.106714473404
$ws.Cells.Item(6, 19).Value = 0.002656813096709221
$ws.Cells.Item(6, 20).Value = 0.002656813096709221

# Row 7
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "Cnr1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 136.676337
$ws.Cells.Item(7, 8).Value = 410.029011
$ws.Cells.Item(7, 9).Value = 0.1929387280825172
$ws.Cells.Item(7, 10).Value = 0.1929387280825172
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.859180333333333
$ws.Cells.Item(7, 14).Value = 8.577541
$ws.Cells.Item(7, 15).Value = 0.9862297573788662
$ws.Cells.Item(7, 16).Value = 0.9862297573788662
$ws.Cells.Item(7, 17).Value = 390.782294782439
$ws.Cells.Item(7, 18).Value = 3517.040653041951
$ws.Cells.Item(7, 19).Value = 0.190281914985808
$ws.Cells.Item(7, 20).Value = 0.190281914985808

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "Cnr1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 143.4557186666667
$ws.Cells.Item(8, 8).Value = 430.367156
$ws.Cells.Item(8, 9).Value = 0.2025088212285795
$ws.Cells.Item(8, 10).Value = 0.2025088212285795
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.03992133333333333
$ws.Cells.Item(8, 14).Value = 0.119764
$ws.Cells.Item(8, 15).Value = 0.01377024262113379
$ws.Cells.Item(8, 16).Value = 0.01377024262113379
$ws.Cells.Item(8, 17).Value = 5.726943563464888
$ws.Cells.Item(8, 18).Value = 51.54249207118399
$ws.Cells.Item(8, 19).Value = 0.002788595601237348
$ws.Cells.Item(8, 20).Value = 0.002788595601237348

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "Cnr1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 143.4557186666667
$ws.Cells.Item(9, 8).Value = 430.367156
$ws.Cells.Item(9, 9).Value = 0.2025088212285795
$ws.Cells.Item(9, 10).Value = 0.2025088212285795
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.859180333333333
$ws.Cells.Item(9, 14).Value = 8.577541
$ws.Cells.Item(9, 15).Value = 0.9862297573788662
$ws.Cells.Item(9, 16).Value = 0.9862297573788662
$ws.Cells.Item(9, 17).Value = 410.1657695159328
$ws.Cells.Item(9, 18).Value = 3691.491925643396
$ws.Cells.Item(9, 19).Value = 0.1997202256273422
$ws.Cells.Item(9, 20).Value = 0.1997202256273422

# Row 10
$ws.Cells.Item(10, 1).Value = "Neutro"
$ws.Cells.Item(10, 2).Value = "Gnai2"
$ws.Cells.Item(10, 3).Value = "Cnr1"
$ws.Cells.Item(10, 4).Value = "FAPs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 198.5982106666667
$ws.Cells.Item(10, 8).Value = 595.794632
$ws.Cells.Item(10, 9).Value = 0.2803505493821544
$ws.Cells.Item(10, 10).Value = 0.2803505493821544
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.03992133333333333
$ws.Cells.Item(10, 14).Value = 0.119764
$ws.Cells.Item(10, 15).Value = 0.01377024262113379
$ws.Cells.Item(10, 16).Value = 0.01377024262113379
$ws.Cells.Item(10, 17).Value = 7.928305367427554
$ws.Cells.Item(10, 18).Value = 71.354748306848
$ws.Cells.Item(10, 19).Value = 0.003860495083960414
$ws.Cells.Item(10, 20).Value = 0.003860495083960414

# Row 11
$ws.Cells.Item(11, 1).Value = "Neutro"
$ws.Cells.Item(11, 2).Value = "Gnai2"
$ws.Cells.Item(11, 3).Value = "Cnr1"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 198.5982106666667
$ws.Cells.Item(11, 8).Value = 595.794632
$ws.Cells.Item(11, 9).Value = 0.2803505493821544
$ws.Cells.Item(11, 10).Value = 0.2803505493821544
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.859180333333333
$ws.Cells.Item(11, 14).Value = 8.577541
$ws.Cells.Item(11, 15).Value = 0.9862297573788662
$ws.Cells.Item(11, 16).Value = 0.9862297573788662
$ws.Cells.Item(11, 17).Value = 567.8280981733235
$ws.Cells.Item(11, 18).Value = 5110.452883559912
$ws.Cells.Item(11, 19).Value = 0.276490054298194
$ws.Cells.Item(11, 20).Value = 0.276490054298194

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Gnai2"
$ws.Cells.Item(12, 3).Value = "Cnr1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 45.55066433333334
$ws.Cells.Item(12, 8).Value = 136.651993
$ws.Cells.Item(12, 9).Value = 0.06430145431675577
$ws.Cells.Item(12, 10).Value = 0.06430145431675577
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.03992133333333333
$ws.Cells.Item(12, 14).Value = 0.119764
$ws.Cells.Item(12, 15).Value = 0.01377024262113379
$ws.Cells.Item(12, 16).Value = 0.01377024262113379
$ws.Cells.Item(12, 17).Value = 1.818443254405778
$ws.Cells.Item(12, 18).Value = 16.365989289652
$ws.Cells.Item(12, 19).Value = 0.0008854466268334774
$ws.Cells.Item(12, 20).Value = 0.0008854466268334774

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Gnai2"
$ws.Cells.Item(13, 3).Value = "Cnr1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 45.55066433333334
$ws.Cells.Item(13, 8).Value = 136.651993
$ws.Cells.Item(13, 9).Value = 0.06430145431675577
$ws.Cells.Item(13, 10).Value = 0.06430145431675577
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.859180333333333
$ws.Cells.Item(13, 14).Value = 8.577541
$ws.Cells.Item(13, 15).Value = 0.9862297573788662
$ws.Cells.Item(13, 16).Value = 0.9862297573788662
$ws.Cells.Item(13, 17).Value = 130.2375636321348
$ws.Cells.Item(13, 18).Value = 1172.138072689213
$ws.Cells.Item(13, 19).Value = 0.0634160076899223
$ws.Cells.Item(13, 20).Value = 0.0634160076899223
